$d = $word.ActiveDocument

# 1) Permission number (MERGEFIELD "Nº" result) in the title paragraph: "2" -> "1"
#    Scope the search to the heading area only, since a bare "2" also occurs
#    later in the document (as part of "02" in a table).
$titlePara = $d.Paragraphs(2).Range
$titlePara.Find.ClearFormatting()
$titlePara.Find.Replacement.ClearFormatting()
$titlePara.Find.Execute("2", $true, $false, $false, $false, $false, $true, 1, $false, "1", 2) | Out-Null

# 2) DNI field result: "45888555" -> "45"
$d.Content.Find.Execute("45888555", $true, $false, $false, $false, $false, $true, 1, $false, "45", 2) | Out-Null

# 3) MODALIDAD field result: "MMO" -> "TEM"
$d.Content.Find.Execute("MMO", $true, $false, $false, $false, $false, $true, 1, $false, "TEM", 2) | Out-Null

# 4) Table cell (ESPACIO CURRICULAR): "EDUC. ARTÍSTICA: MÚSICA" -> "PASANTÍA"
$d.Content.Find.Execute("EDUC. ARTÍSTICA: MÚSICA", $true, $false, $false, $false, $false, $true, 1, $false, "PASANTÍA", 2) | Out-Null

# 5) Table cell (CURSO): "1°2°" -> "7°2°"
$d.Content.Find.Execute("1°2°", $true, $false, $false, $false, $false, $true, 1, $false, "7°2°", 2) | Out-Null

# 6) Closing date: "Tinogasta, 02 de diciembre de 2024" -> "Tinogasta, 03 de diciembre de 2024"
$d.Content.Find.Execute("Tinogasta, 02 de diciembre de 2024", $true, $false, $false, $false, $false, $true, 1, $false, "Tinogasta, 03 de diciembre de 2024", 2) | Out-Null
